$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture a source format (fontId=2 "Arial 10 regular no-border/fill" look) ---
# B4 in the original sheet ("Publicly available report or research article") and D4
# both carry that look; D4 is not touched by this edit, so grab the format from there
# before we start rewriting column B.
$ws.Range("D4").Copy()
$ws.Range("B6").PasteSpecial(-4122)   # xlPasteFormats

# --- Rewrite the lookup values in column B (dataset rows split into static/dynamic) ---
$ws.Range("B2").Value = "Publicly available dataset, static web page"
$ws.Range("B3").Value = "Publicly available dataset, dynamic web page"
$ws.Range("B4").Value = "Proprietary dataset, static web page"
$ws.Range("B5").Value = "Proprietary dataset, dynamic web page"
$ws.Range("B6").Value = "Publicly available report or research article"

# New row 6 (sheet row 7): shift the former "Proprietary report or research article"
# entry down and give it id 6.
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Proprietary report or research article"

# --- Formatting touch-ups ---
# B4 / B5 should look like B2 (Arial 10 regular, same as the rest of column B).
$ws.Range("B2").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B5").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Rename the default cell style, mirroring the source workbook's "Standard" name.
$wb.Styles.Item(1).Name = "Standard"

# --- Selection state recorded by the author's last save ---
$ws.Range("A2:B7").Select() | Out-Null
